$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column G (SMD), shifting old G (effectSize) -> H and old H (significance) -> I
$ws.Columns("G:G").Insert()

# Approximate the target display width for new column G (target xml width 12.7109375)
$ws.Columns("G:G").ColumnWidth = 11.9

# Header
$ws.Range("G1").Value = "SMD"

# Row 2
$ws.Range("G2").Value = 55.39165824783241
$ws.Range("D2").Value = 375858.88592197088
$ws.Range("F2").Value = 0.99869801660552382

# Row 3
$ws.Range("G3").Value = 0.18090670659715877
$ws.Range("D3").Value = 4.0090864702492341
$ws.Range("E3").Value = 0.045806251947808296
$ws.Range("F3").Value = 0.0081154103842393053

# Row 4
$ws.Range("G4").Value = 2.7787675242321499
$ws.Range("D4").Value = 945.88974683159074
$ws.Range("F4").Value = 0.65874817263565988

# Row 5
$ws.Range("G5").Value = 1.0648117432538804
$ws.Range("D5").Value = 138.89344594999253
$ws.Range("F5").Value = 0.22085370239497912
